$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/J1 should carry the same (bold/bordered) style as the
# existing header row. Copy H1's formatting into I1:J1 first, then overwrite
# the values, so the destination cells keep the "real" header look/style
# index rather than getting a brand-new (if equivalent) style slot.
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-42: I = constant 1, J = mirrors column H (IP) for that row -
# except row 3, which the source data records as I3=9, J3=9.
for ($r = 2; $r -le 42; $r++) {
    if ($r -eq 3) {
        $ws.Cells.Item($r, 9).Value = 9
        $ws.Cells.Item($r, 10).Value = 9
    } else {
        $hVal = $ws.Cells.Item($r, 8).Value()
        $ws.Cells.Item($r, 9).Value = 1
        $ws.Cells.Item($r, 10).Value = $hVal
    }
}
